# Update "想去人数" (number of attendees) figures that changed between scrapes.
# Sheet "展览" (Exhibitions)
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 100
$ws1.Range("F3").Value = 174
$ws1.Range("F4").Value = 410
$ws1.Range("F5").Value = 185
$ws1.Range("F6").Value = 128
$ws1.Range("F7").Value = 1093
$ws1.Range("F8").Value = 366
$ws1.Range("F9").Value = 189
$ws1.Range("F13").Value = 369
$ws1.Range("F15").Value = 157
$ws1.Range("F16").Value = 716
$ws1.Range("F18").Value = 72
$ws1.Range("F19").Value = 991
$ws1.Range("F20").Value = 447
$ws1.Range("F21").Value = 256
$ws1.Range("F26").Value = 462

# Sheet "本地生活" (Local life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 343

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 343
$ws4.Range("F4").Value = 100
$ws4.Range("F5").Value = 174
$ws4.Range("F6").Value = 410
$ws4.Range("F7").Value = 185
$ws4.Range("F8").Value = 128
$ws4.Range("F9").Value = 1093
$ws4.Range("F10").Value = 366
$ws4.Range("F11").Value = 189
$ws4.Range("F20").Value = 369
$ws4.Range("F22").Value = 157
$ws4.Range("F23").Value = 716
$ws4.Range("F25").Value = 72
$ws4.Range("F26").Value = 991
$ws4.Range("F27").Value = 447
$ws4.Range("F30").Value = 256
$ws4.Range("F33").Value = 627
$ws4.Range("F38").Value = 462
